$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.933.04"
$ws.Range("E2").Value = "  +5.72%  "
$ws.Range("D3").Value = "2.522.90"
$ws.Range("E3").Value = "  +9.02%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'485.67"
$ws.Range("E5").Value = "  +11.37%  "
$ws.Range("D6").Value = "'142.82"
$ws.Range("E6").Value = "  +18.76%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.515"
$ws.Range("E8").Value = "  +8.98%  "
$ws.Range("D9").Value = "2.521.09"
$ws.Range("E9").Value = "  +8.84%  "
$ws.Range("D10").Value = "'0.0997"
$ws.Range("E10").Value = "  +10.19%  "
$ws.Range("D11").Value = "'5.53"
$ws.Range("E11").Value = "  +5.34%  "
$ws.Range("E12").Value = "  +8.14%  "
$ws.Range("D13").Value = "'0.123"
$ws.Range("E13").Value = "  +0.97%  "
$ws.Range("D14").Value = "2.954.74"
$ws.Range("E14").Value = "  +9.07%  "
$ws.Range("D15").Value = "55.940.76"
$ws.Range("E15").Value = "  +5.75%  "
$ws.Range("D16").Value = "'20.77"
$ws.Range("E16").Value = "  +9.32%  "
$ws.Range("E17").Value = "  +17.03%  "
$ws.Range("D18").Value = "2.516.63"
$ws.Range("E18").Value = "  +7.34%  "
$ws.Range("D19").Value = "'4.42"
$ws.Range("E19").Value = "  +11.99%  "
$ws.Range("D20").Value = "'322.49"
$ws.Range("E20").Value = "  +7.62%  "
$ws.Range("D21").Value = "'10.08"
$ws.Range("E21").Value = "  +11.33%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").Value = "'5.76"
$ws.Range("E23").Value = "  +7.85%  "
$ws.Range("D24").Value = "'58.34"
$ws.Range("E24").Value = "  +5.48%  "
$ws.Range("E25").Value = "  +10.03%  "
$ws.Range("D26").Value = "'0.411"
$ws.Range("E26").Value = "  +12.19%  "
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").Value = "2.615.60"
$ws.Range("E28").Value = "  +7.83%  "
$ws.Range("E29").Value = "  +8.45%  "
$ws.Range("D30").Value = "0.0₃0813"
$ws.Range("E30").Value = "  +16.68%  "
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("D32").Value = "'150.14"
$ws.Range("E32").Value = "  +5.33%  "
$ws.Range("D33").Value = "'18.30"
$ws.Range("E33").Value = "  +7.06%  "
$ws.Range("D34").Value = "'1.50"
$ws.Range("E34").Value = "  +11.52%  "
$ws.Range("D35").Value = "'5.24"
$ws.Range("E35").Value = "  +10.76%  "
$ws.Range("D36").Value = "'0.888"
$ws.Range("E36").Value = "  +7.50%  "
$ws.Range("D37").Value = "'3.75"
$ws.Range("E37").Value = "  +6.59%  "
$ws.Range("E38").Value = "  +13.97%  "
$ws.Range("D39").Value = "'34.41"
$ws.Range("E39").Value = "  +3.28%  "
$ws.Range("D40").Value = "'0.619"
$ws.Range("E40").Value = "  +18.76%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.48%  "
$ws.Range("D42").Value = "'0.0558"
$ws.Range("E42").Value = "  +12.07%  "
$ws.Range("D43").Value = "'3.45"
$ws.Range("E43").Value = "  +9.46%  "
$ws.Range("E44").Value = "  +10.82%  "
$ws.Range("D45").Value = "2.010.65"
$ws.Range("E45").Value = "  +6.61%  "
$ws.Range("D46").Value = "'4.68"
$ws.Range("E46").Value = "  +21.89%  "
$ws.Range("D47").Value = "'10.15"
$ws.Range("E47").Value = "  -1.19%  "
$ws.Range("D48").Value = "'0.0910"
$ws.Range("E48").Value = "  +9.45%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0226"
$ws.Range("E49").Value = "  +8.70%  "
$ws.Range("B50").Value = "Bittensor"
$ws.Range("C50").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D50").Value = "'254.47"
$ws.Range("E50").Value = "  +36.04%  "
$ws.Range("D51").Value = "'17.81"
$ws.Range("E51").Value = "  +14.78%  "
